try {
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3: new publication entry (JAMA Network Open / SNAP emergency allotments) ----

# A3: authors, rich text with per-run coloring matching the site's markdown author byline style
$authorsText = "**Matthew Lavallee**, Sandro Galea, Nadia N. Abuelezam"
$ws.Range("A3").Value = $authorsText
$ws.Range("A3").Font.Color = 4473924
$ws.Range("A3").Characters(22,1).Font.Color = 3355443
$ws.Range("A3").Characters(23,13).Font.Color = 4473924
$ws.Range("A3").Characters(36,1).Font.Color = 3355443
$ws.Range("A3").Characters(37,18).Font.Color = 4473924

# B3: year
$ws.Range("B3").Value = 2023

# C3: title
$ws.Range("C3").Value = "Supplemental Nutrition Assistance Program Emergency Allotments and Food Security, Hospitalizations, and Hospital Capacity"

# D3: journal
$ws.Range("D3").Value = "*JAMA Network Open*"

# F3: publication date (2023-08-12)
$ws.Range("F3").Value = 45150

# G3: link
$ws.Range("G3").Value = "https://jamanetwork.com/journals/jamanetworkopen/fullarticle/2808124"

# H3: short_id
$ws.Range("H3").Value = "b"

# E3: abstract (populated last to match original shared-string ordering)
$ws.Range("E3").Value = "`n**IMPORTANCE:** Understanding how social policies shape health is a national priority, especially in the context of the COVID-19 pandemic.`n**OBJECTIVE:** To understand the association between politically motivated changes to Nebraska’s Supplemental Nutrition Assistance Program (SNAP) policy and public health measures during the COVID-19 pandemic.`n**DESIGN, STUDY, AND PARTICIPANTS:** This cross-sectional study used synthetic control methods to estimate the association of Nebraska’s decision to reject emergency allotments for the SNAP with food security and hospital capacity indicators. A counterfactual for Nebraska was created by weighting data from the rest of the US. State-level changes in Nebraska between March 2020 and March 2021 were included. Data were acquired from the Census Bureau’s Household Pulse Survey on individual food security and mental health indicators and from the US Centers for Disease Control and Prevention on hospital-level capacity indicators. Data analysis occurred between October 2022 and June 2023.`n**INTERVENTION:* The rejection of additional SNAP funds for low-income households in Nebraska from August to November 2020.`n**MAIN OUTCOMES AND MEASURES:** Food insecurity and inpatient bed use indicators (ie, inpatient beds filled, inpatient beds filled by patients with COVID-19, and inpatients with COVID-19).`n**RESULTS:** The survey data of 1 591 006 respondents from May 2020 to November 2020 was analyzed, and 24 869 (1.56%) lived in Nebraska. Nebraska’s population was composed of proportionally more White individuals (mean [SD], 88.70% [0.29%] vs 78.28% [0.26%]; P < .001), fewer individuals who made more than `$200 000 in 2019 (4.20% [0.45%] vs 5.22% [0.12%]; P < .001), and more households sized 1 to 3 (63.41% [2.29%] vs 61.13% [1.10%); P = .03) compared with other states. Nebraska’s rejection of additional funding for SNAP recipients was associated with increases in food insecurity (raw mean [SD] difference 1.61% [1.30%]; relative difference, 19.63%; P = .02), percentage of inpatient beds filled by patients with COVID-19 (raw mean [SD] difference, 0.19% [1.55%]; relative difference, 3.90%; P = .02), and percentage of inpatient beds filled (raw mean [SD] difference, 2.35% [1.82%]; relative difference, 4.10%; P = .02).`n**CONCLUSIONS AND RELEVANCE:** In this cross-sectional study, the association between social policy, food security, health, and public health resources was examined, and the rejection of emergency allotments in Nebraska was associated with increased food insecurity. Additionally, this intervention was associated with an increased rate of hospitalizations for COVID-19 and non–COVID-19 causes."

# Row 3 needs the same max row height as row 2 (long wrapped abstract text)
$ws.Rows.Item(3).RowHeight = 409.6

# Update the active selection to the new row, like Excel would after entering data
$ws.Range("E3").Select()

} catch {
    Write-Host "ERROR:" $_.Exception.Message
}
